# Realestate Update resale numbers 2025-02-23 16:46
# Appends a new data row (row 84) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 84

# Text columns: leading apostrophe forces these to stay as text
# (otherwise "2025-02-23" would become a date serial and "08" would
# become the number 8, dropping the leading zero). Resetting the
# Style back to "Normal" after the assignment avoids leaving an
# explicit text-number-format style on the cell.
$ws.Cells.Item($row, 1).Value = "'2025-02-23"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "16:46:51"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "'08"
$ws.Cells.Item($row, 4).Style = "Normal"

# Numeric columns
$ws.Cells.Item($row, 5).Value = 130667
$ws.Cells.Item($row, 6).Value = 141856
$ws.Cells.Item($row, 7).Value = 172561
$ws.Cells.Item($row, 8).Value = 158716
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146655
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193566
$ws.Cells.Item($row, 14).Value = 115508
$ws.Cells.Item($row, 15).Value = 46449
$ws.Cells.Item($row, 16).Value = 29340
$ws.Cells.Item($row, 17).Value = 68540
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48854
$ws.Cells.Item($row, 20).Value = -1
